# Update countries & provincias Spain
# Applies the "5 de Abril de 2020" data refresh: updates the timestamp caption,
# refreshes several countries' stats, and inserts Serbia's updated figures in
# its new sorted position (ahead of Mexico), shifting Mexico/Panama/Peru/Grecia
# down one row while Sudafrica's row remains unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header timestamp -------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 5 de Abril de 2020 a las 15:52"

# --- Estados Unidos (row 4) --------------------------------------------
$ws.Range("B4").Value = 311656
$ws.Range("C4").Value = 299
$ws.Range("E4").Value = 288374

# --- Alemania (row 7) ----------------------------------------------------
$ws.Range("B7").Value = 97074
$ws.Range("C7").Value = 982
$ws.Range("E7").Value = 69196

# --- Reino Unido (row 11) -------------------------------------------------
$ws.Range("B11").Value = 47806
$ws.Range("C11").Value = 5903
$ws.Range("E11").Value = 42739
$ws.Range("G11").Value = 619
$ws.Range("H11").Value = 4932

# --- Suiza (row 13) --------------------------------------------------------
$ws.Range("E13").Value = 14000
$ws.Range("G13").Value = 19
$ws.Range("H13").Value = 685

# --- Suecia (row 22) --------------------------------------------------------
$ws.Range("F22").Value = 406

# --- Rows 43-47: Serbia moves ahead of Mexico with refreshed numbers,
#     Mexico/Panama/Peru/Grecia shift down a row, Sudafrica (row 48) unchanged.
$ws.Range("A43").Value = "Serbia"
$ws.Range("B43").Value = 1908
$ws.Range("C43").Value = 284
$ws.Range("D43").Value = 54
$ws.Range("E43").Value = 1803
$ws.Range("F43").Value = 98
$ws.Range("G43").Value = 7
$ws.Range("H43").Value = 51

$ws.Range("A44").Value = "Mexico"
$ws.Range("B44").Value = 1890
$ws.Range("C44").Value = 202
$ws.Range("D44").Value = 633
$ws.Range("E44").Value = 1178
$ws.Range("F44").Value = 1
$ws.Range("G44").Value = 19
$ws.Range("H44").Value = 79

$ws.Range("A45").Value = "Panama"
$ws.Range("B45").Value = 1801
$ws.Range("C45").Value = 0
$ws.Range("D45").Value = 13
$ws.Range("E45").Value = 1742
$ws.Range("F45").Value = 75
$ws.Range("G45").Value = 0
$ws.Range("H45").Value = 46

$ws.Range("A46").Value = "Peru"
$ws.Range("B46").Value = 1746
$ws.Range("C46").Value = 0
$ws.Range("D46").Value = 914
$ws.Range("E46").Value = 759
$ws.Range("F46").Value = 88
$ws.Range("G46").Value = 0
$ws.Range("H46").Value = 73

$ws.Range("A47").Value = "Grecia"
$ws.Range("B47").Value = 1673
$ws.Range("C47").Value = 0
$ws.Range("D47").Value = 78
$ws.Range("E47").Value = 1525
$ws.Range("F47").Value = 92
$ws.Range("G47").Value = 2
$ws.Range("H47").Value = 70

# --- Croacia (row 58) -----------------------------------------------------
$ws.Range("E58").Value = 1048
$ws.Range("G58").Value = 3
$ws.Range("H58").Value = 15

# --- Sri Lanka (row 109) ---------------------------------------------------
$ws.Range("B109").Value = 174
$ws.Range("C109").Value = 8
$ws.Range("E109").Value = 140
